$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.099.35'

$ws.Range("D3").Value = '2.978.99'
$ws.Range("E3").Value = '  -0.48%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.97'
$ws.Range("E5").Value = '  +1.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.69'
$ws.Range("E6").Value = '  -1.81%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.514'
$ws.Range("E8").Value = '  -1.27%  '

$ws.Range("D9").Value = '2.976.20'
$ws.Range("E9").Value = '  -0.60%  '

$ws.Range("E10").Value = '  -0.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.00'
$ws.Range("E11").Value = '  +5.34%  '

$ws.Range("E12").Value = '  +2.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000225'
$ws.Range("E13").Value = '  -0.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.07'
$ws.Range("E14").Value = '  -1.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.126'
$ws.Range("E15").Value = '  +2.42%  '

$ws.Range("D16").Value = '3.471.00'
$ws.Range("E16").Value = '  -0.50%  '

$ws.Range("D17").Value = '61.233.96'
$ws.Range("E17").Value = '  -1.67%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.84'
$ws.Range("E18").Value = '  -1.67%  '

$ws.Range("D19").Value = '2.977.36'
$ws.Range("E19").Value = '  -0.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.15'
$ws.Range("E20").Value = '  -1.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.99'
$ws.Range("E21").Value = '  +1.39%  '

$ws.Range("E22").Value = '  +0.53%  '

$ws.Range("E23").Value = '  -1.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.82'
$ws.Range("E24").Value = '  +2.00%  '

$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.43'
$ws.Range("E25").Value = '  +2.60%  '

$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.16'
$ws.Range("E26").Value = '  -3.96%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.89'
$ws.Range("E27").Value = '  -2.45%  '

$ws.Range("E28").Value = '  +0.02%  '

$ws.Range("E29").Value = '  +2.48%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.09%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.12'
$ws.Range("E31").Value = '  -0.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.04'
$ws.Range("E32").Value = '  -1.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.18'
$ws.Range("E33").Value = '  +1.60%  '

$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").Value = '0.0₃0805'
$ws.Range("E35").Value = '  +2.85%  '

$ws.Range("E36").Value = '  -0.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.76'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.00'
$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("E39").Value = '  -2.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.93'
$ws.Range("E40").Value = '  +0.19%  '

$ws.Range("E41").Value = '  +6.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.81'
$ws.Range("E42").Value = '  -2.97%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '384.24'
$ws.Range("E43").Value = '  -2.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.267'
$ws.Range("E44").Value = '  -2.23%  '

$ws.Range("E45").Value = '  -0.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.20'
$ws.Range("E46").Value = '  -1.70%  '

$ws.Range("D47").Value = '2.695.14'
$ws.Range("E47").Value = '  -2.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.45'
$ws.Range("E48").Value = '  +1.65%  '

$ws.Range("E50").Value = '  -0.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.12'
$ws.Range("E51").Value = '  -0.54%  '
